$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (values + copy the header formatting from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row data: row, I-value, J-value
$data = @(
    @(2,1,5),
    @(3,10,13),
    @(4,3,7),
    @(5,9,9),
    @(6,1,5),
    @(7,1,6),
    @(8,1,4),
    @(9,1,6),
    @(10,1,6),
    @(11,1,4),
    @(12,1,2),
    @(13,1,6),
    @(14,1,6),
    @(15,1,5),
    @(16,1,5),
    @(17,1,7),
    @(18,1,4),
    @(19,1,7),
    @(20,1,5),
    @(21,1,7),
    @(22,1,5),
    @(23,1,5),
    @(24,1,5),
    @(25,1,6),
    @(26,1,6),
    @(27,1,5),
    @(28,1,7),
    @(29,1,6),
    @(30,1,5),
    @(31,1,5),
    @(32,1,6),
    @(33,1,4),
    @(34,1,6),
    @(35,1,5),
    @(36,1,4),
    @(37,1,4),
    @(38,1,3),
    @(39,1,2)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
